$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New quest line rows appended to the Items sheet (quests that can take copper coins)
$rows = @(
    @{ Row = 72; D = "Musty old Church Records"; E = "quest"; F = "These old Church records state that alchemy was apart of the Church in one way or another. It seems, if I can make out these musty old pages, that the Church tried to use Holy magic and Alchemy to create special types of candles."; L = 0; M = 0; V = 0; W = 0; BD = "Dungeons of Valifore"; BE = 0.05; BF = 1 },
    @{ Row = 73; D = "Alchemically Enchanted Holy Candle"; E = "quest"; F = "This item will let you ignore the Leveling caps and gain XP faster."; L = 0; M = 0; V = 0; W = 0; BE = 0.95; BF = 1 },
    @{ Row = 74; D = "Dried up Enchanted Holy Ink"; E = "quest"; F = "This old ink was used to write enchantments on candles back when the Church was still a thing through out Tlessa, fanatical bunch if you ask me. They created these types of enchanted `"Holy Inks`" that only the elite could get their hands on. The Candle Maker might be able to make use of this."; L = 0; M = 0; V = 0; W = 0; BD = "Shadow Caves"; BE = 0.15; BF = 1 },
    @{ Row = 75; D = "Fanatics Candle of Despair"; E = "quest"; F = "Well, I guess that ink you had found was a bit cursed. But this should still do the trick child. I think ..."; L = 0; M = 0; V = 0; W = 0; BE = 1.25; BF = 1 },
    @{ Row = 76; D = "Corrupted Candle of the Church"; E = "quest"; F = "Well it's become corrupted, the alchemical process failed, but it seems to be emanating a strange glow of power."; L = 0; M = 0; V = 0; W = 0; BE = 1.6; BF = 1 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 22).Value = $r.V
    $ws.Cells.Item($row, 23).Value = $r.W
    if ($r.ContainsKey("BD")) {
        $ws.Cells.Item($row, 56).Value = $r.BD
    }
    $ws.Cells.Item($row, 57).Value = $r.BE
    $ws.Cells.Item($row, 58).Value = $r.BF
}

# Column D needs to widen to fit the new, longer quest-item names
# (ColumnWidth is expressed in character units that get padding of 5/6
# added on save, so subtract that back off to land on an exact 41)
$ws.Columns.Item(4).ColumnWidth = 41 - 5/6
